$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.166.22'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.54%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.902.09'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.10%  '
$ws.Range("E4").Value = '  +0.30%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '305.97'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.25%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.29%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5259'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.67%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3776'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.55%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07248'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.73%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.17'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.28%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8984'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.39%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08316'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +9.90%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.906.08'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.49%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '94.77'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.43%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.267'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.11%  '
$ws.Range("E17").Value = '  +1.36%  '
$ws.Range("E18").Value = '  +1.60%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '27.202.07'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.53%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.059'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.52%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.130.90'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.02%  '
$ws.Range("E23").Value = '  +1.48%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.426'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.38%  '
$ws.Range("B25").Value = 'LidoDAOToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.279'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +7.35%  '
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '146.40'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.43%  '
$ws.Range("E27").Value = '  -1.31%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.11'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.56%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '114.75'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.11%  '
$ws.Range("E30").Value = '  -0.04%  '
$ws.Range("E31").Value = '  -0.18%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09258'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.63%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.8117'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +6.57%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05052'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.38%  '
$ws.Range("E35").Value = '  +4.05%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.980'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.64%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.331'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.71%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.582'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.40%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5710'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.90%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01979'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.65%  '
$ws.Range("E41").Value = '  -0.14%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.660'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.13%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.950'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.97%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '118.20'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.48%  '
$ws.Range("E45").Value = '  +0.56%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4834'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.97%  '
$ws.Range("E47").Value = '  +0.38%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.14'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.00%  '
$ws.Range("E49").Value = '  +2.76%  '
$ws.Range("E50").Value = '  +1.00%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '63.53'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.14%  '
